# Update Hgf-Sdc1.xlsx with the re-run ("new TPM") NATMI numbers.
#
# The old sheet had 9 data rows (3 sending clusters x 3 target clusters:
# ECs, FAPs, MuSCs all sending). The new run drops "ECs" as a sending
# cluster, leaving only FAPs and MuSCs as senders (2 x 3 = 6 data rows),
# and refreshes every numeric column with the recomputed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the three rows whose "Sending cluster" (column A) was ECs - this
# shifts the old FAPs-sending / MuSCs-sending rows up to rows 2-7 and
# shrinks the used range from A1:T10 to A1:T7.
$ws.Range("A2:T4").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp)

# Each inner array is one data row, in column order:
# A, D, E, F, G, H, I, J, K, L, M, N, O, P, Q, R, S, T
# (B is always "Hgf", C is always "Sdc1" - unchanged from before)
$newData = @(
  @(2, "FAPs",  "ECs",   3, 1, "8.583520999999999",  "25.750563",  "0.8910607110509009", "0.8910607110509009", 3, 1, "0.8213140000000001", "2.463942",  "0.06824749762056036", "0.06824749762056037", "7.049765966594",   "63.447893699346",  "0.06081266375722118", "0.06081266375722119"),
  @(3, "FAPs",  "FAPs",  3, 1, "8.583520999999999",  "25.750563",  "0.8910607110509009", "0.8910607110509009", 3, 1, "2.427350333333333",  "7.282051",  "0.2017018900182306",  "0.2017018900182306",  "20.83521256052367", "187.516913044713", "0.1797286295399552",  "0.1797286295399552"),
  @(4, "FAPs",  "MuSCs", 3, 1, "8.583520999999999",  "25.750563",  "0.8910607110509009", "0.8910607110509009", 3, 1, "8.785681666666667",  "26.357045", "0.7300506123612091",  "0.7300506123612091",  "75.41208308514832", "678.708747766335", "0.6505194177537246",  "0.6505194177537246"),
  @(5, "MuSCs", "ECs",   3, 1, "1.049404",            "3.148212",   "0.1089392889490991", "0.1089392889490991", 3, 1, "0.8213140000000001", "2.463942",  "0.06824749762056036", "0.06824749762056037", "0.8618901968560001","7.757011771704001","0.007434833863339175","0.007434833863339177"),
  @(6, "MuSCs", "FAPs",  3, 1, "1.049404",            "3.148212",   "0.1089392889490991", "0.1089392889490991", 3, 1, "2.427350333333333",  "7.282051",  "0.2017018900182306",  "0.2017018900182306",  "2.547271149201333", "22.925440342812",  "0.02197326047827542", "0.02197326047827542"),
  @(7, "MuSCs", "MuSCs", 3, 1, "1.049404",            "3.148212",   "0.1089392889490991", "0.1089392889490991", 3, 1, "8.785681666666667",  "26.357045", "0.7300506123612091",  "0.7300506123612091",  "9.219729483726667", "82.97756535354",   "0.07953119460748447", "0.07953119460748447")
)

foreach ($d in $newData) {
  $rowNum = $d[0]
  $ws.Cells.Item($rowNum, 1).Value  = $d[1]            # A: Sending cluster
  $ws.Cells.Item($rowNum, 2).Value  = "Hgf"             # B: Ligand symbol (unchanged)
  $ws.Cells.Item($rowNum, 3).Value  = "Sdc1"            # C: Receptor symbol (unchanged)
  $ws.Cells.Item($rowNum, 4).Value  = $d[2]             # D: Target cluster
  $ws.Cells.Item($rowNum, 5).Value  = [double]$d[3]     # E
  $ws.Cells.Item($rowNum, 6).Value  = [double]$d[4]     # F
  $ws.Cells.Item($rowNum, 7).Value  = [double]$d[5]     # G
  $ws.Cells.Item($rowNum, 8).Value  = [double]$d[6]     # H
  $ws.Cells.Item($rowNum, 9).Value  = [double]$d[7]     # I
  $ws.Cells.Item($rowNum, 10).Value = [double]$d[8]     # J
  $ws.Cells.Item($rowNum, 11).Value = [double]$d[9]     # K
  $ws.Cells.Item($rowNum, 12).Value = [double]$d[10]    # L
  $ws.Cells.Item($rowNum, 13).Value = [double]$d[11]    # M
  $ws.Cells.Item($rowNum, 14).Value = [double]$d[12]    # N
  $ws.Cells.Item($rowNum, 15).Value = [double]$d[13]    # O
  $ws.Cells.Item($rowNum, 16).Value = [double]$d[14]    # P
  $ws.Cells.Item($rowNum, 17).Value = [double]$d[15]    # Q
  $ws.Cells.Item($rowNum, 18).Value = [double]$d[16]    # R
  $ws.Cells.Item($rowNum, 19).Value = [double]$d[17]    # S
  $ws.Cells.Item($rowNum, 20).Value = [double]$d[18]    # T
}
